# Update the "Overview" income-statement sheet with the actual reported
# figures (replacing the placeholder zeros / dash markers) for the
# methanol/shekhark yearly income statement (rial).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row -> values for columns D,E,F,G,H
# (12-month periods ending 1397/12, 1398/12, 1399/12, 1400/12, 1401/12)
$data = @{
    11 = @(29622110, 26193145, 51706449, 93012073, 111798328)      # فروش
    12 = @(-9722364, -11291404, -18569571, -55073067, -66958743)   # بهای تمام شده کالای فروش رفته
    13 = @(19899746, 14901741, 33136878, 37939006, 44839585)       # سود (زیان) ناخالص
    14 = @(-2470746, -3583454, -9611879, -9757345, -10713596)      # هزینه های عمومی, اداری و تشکیلاتی
    15 = @(-677609, 0, 0, 0, 0)                                    # هزینه کاهش ارزش دریافتنی‌ها
    16 = @(7904263, 6812396, 8661548, 1643256, 10344702)           # خالص سایر درامدها (هزینه ها) ی عملیاتی
    17 = @(24655654, 18130683, 32186547, 29824917, 44470691)       # سود (زیان) عملیاتی
    18 = @(0, 0, 0, 0, 0)                                          # هزینه های مالی
    19 = @(832533, 3319261, 1663178, 1911565, 2137418)             # خالص سایر درامدها و هزینه های غیرعملیاتی
    20 = @(25488187, 21449944, 33849725, 31736482, 46608109)       # سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
    21 = @(0, -7191, -30982, -3242156, -8134088)                   # مالیات
    22 = @(25488187, 21442753, 33818743, 28494326, 38474021)       # سود (زیان) خالص عملیات در حال تداوم
    23 = @(0, 0, 0, 0, 0)                                          # سود (زیان) عملیات متوقف شده پس از اثر مالیاتی
    24 = @(25488187, 21442753, 33818743, 28494326, 38474021)       # سود (زیان) خالص
    25 = @(4248, 3574, 5636, 4749, 6412)                           # سود هر سهم پس از کسر مالیات
    26 = @(6000000, 6000000, 6000000, 6000000, 6000000)            # سرمایه
    27 = @(4248, 3574, 5636, 4749, 6412)                           # سود هر سهم بر اساس آخرین سرمایه
}

$cols = @("D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
